# Weekly data refresh: a new weekly reading is inserted as the new row 3
# (the rest of the historical readings shift down by one row, row 2 is
# untouched), matching the "Fruta / hortaliza, semanal" update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing rows 3:28 down to 4:29, inserting a blank row 3
# (formatting of the row above is carried down automatically, same as
# Excel's native Insert behaviour).
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with this week's reading.
$ws.Cells.Item(3, 1).Value2  = 11
$ws.Cells.Item(3, 2).Value2  = "Vega Monumental Concepción"
$ws.Cells.Item(3, 3).Value2  = "Bíobío"
$ws.Cells.Item(3, 4).Value2  = 44691
$ws.Cells.Item(3, 5).Value2  = 8
$ws.Cells.Item(3, 6).Value2  = 100114007
$ws.Cells.Item(3, 7).Value2  = "Jengibre"
$ws.Cells.Item(3, 8).Value2  = "Sin especificar"
$ws.Cells.Item(3, 9).Value2  = "Primera"
$ws.Cells.Item(3, 10).Value2 = 100
$ws.Cells.Item(3, 11).Value2 = 12000
$ws.Cells.Item(3, 12).Value2 = 13000
$ws.Cells.Item(3, 13).Value2 = 12500
$ws.Cells.Item(3, 14).Value2 = "$/caja 13 kilos"
$ws.Cells.Item(3, 15).Value2 = "Perú"
$ws.Cells.Item(3, 16).Value2 = 962
$ws.Cells.Item(3, 17).Value2 = 13
$ws.Cells.Item(3, 18).Value2 = "Hortaliza"
